$wb = $excel.ActiveWorkbook
$wsUK = $wb.Worksheets.Item("UK")
$wsHungary = $wb.Worksheets.Item("Hungary")

# --- UK sheet: fill in the User Story / Jira reference cell (B4) ---
# Copy formatting from D3 (which already uses the wrap-text style) then set the text.
$wsUK.Range("D3").Copy()
$wsUK.Range("B4").PasteSpecial(-4122)
$wsUK.Range("B4").Value = "NGC-3003/ T1240/T1246/T1255"
$wsUK.Rows("4:4").RowHeight = 43.2

# --- UK sheet: insert a new "IOB800" row above the existing "XIOM" row (row 10) ---
$wsUK.Rows("10:10").Insert()
$wsUK.Range("A9").Copy()
$wsUK.Range("A10").PasteSpecial(-4122)
$wsUK.Range("A10").Value = "IOB800"

# --- Selection / active-tab bookkeeping ---
# Hungary loses the tab-selected / active-cell state it had before...
$wsHungary.Range("A7:A17").Select() | Out-Null
# ...and UK becomes the active sheet/tab with B4 selected.
$wsUK.Activate()
$wsUK.Range("B4").Select() | Out-Null
